$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: merge the two runs "FRI Nov 23" + " 11:52:07 IST 2018" into a
# single run "FRI Nov 23 11:52:07 IST 2018" (same visible text, just a
# run-splitting cleanup). Find/Replace across the run boundary naturally
# collapses it into one run with the paragraph's existing formatting.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "FRI Nov 23 11:52:07 IST 2018", $false, $false, $false, $false, $false,
    $true, 1, $false, "FRI Nov 23 11:52:07 IST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: append a new "30/11/2018 MAMATHA CHICK IN" purchase record
# after the most recent "Amount Received mode ... - CASH AND CLEARD"
# block (the one nearest the end of the document, right before the
# trailing blank paragraphs).
# ---------------------------------------------------------------------

# Locate the last paragraph whose text mentions "CASH AND CLEARD" - this
# is the end of the most recent purchase record.
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*CASH AND CLEARD*") {
        $anchorIdx = $i
    }
}

$global:cur = $anchorIdx

function Add-NewParagraph {
    $d.Paragraphs($global:cur).Range.InsertParagraphAfter()
    $global:cur = $global:cur + 1
}

function Add-Text([string]$text) {
    $d.Paragraphs($global:cur).Range.InsertAfter($text)
}

function Add-Tabs([int]$count) {
    for ($t = 0; $t -lt $count; $t++) {
        $d.Paragraphs($global:cur).Range.InsertAfter([char]9)
    }
}

function Set-Bold {
    $d.Paragraphs($global:cur).Range.Bold = 1
}

function Clear-Bold {
    $d.Paragraphs($global:cur).Range.Bold = 0
}

# Blank separator paragraph
Add-NewParagraph

# Timestamp line
Add-NewParagraph
Add-Text "THU Nov 29 11:44:46 IST 2018"

# Person Name	- CHANDU
Add-NewParagraph
Add-Text "Person Name"
Add-Tabs 4
Add-Text "- CHANDU"

# Bill number	- 9178
Add-NewParagraph
Add-Text "Bill number"
Add-Tabs 4
Add-Text "- 9178"

# separator line
Add-NewParagraph
Add-Text "---------------------------------------------------------------"

# Item Name	- SORE KAI
Add-NewParagraph
Add-Text "Item Name"
Add-Tabs 4
Add-Text "- SORE KAI"

# Number of Pockets	- 2
Add-NewParagraph
Add-Text "Number of Pockets"
Add-Tabs 3
Add-Text "- 2"

# Number of KGs	- 111
Add-NewParagraph
Add-Text "Number of KGs"
Add-Tabs 3
Add-Text "- 111"

# Rate	- 8
Add-NewParagraph
Add-Text "Rate"
Add-Tabs 5
Add-Text "- 8"

# Total Price	- 888.0
Add-NewParagraph
Add-Text "Total Price"
Add-Tabs 4
Add-Text "- 888.0"

# Amount balance	- 888.0 (bold)
Add-NewParagraph
Set-Bold
Add-Text "Amount balance"
Add-Tabs 3
Add-Text "- 888.0"

# Trailing blank paragraph (not bold, even though it follows a bold one)
Add-NewParagraph
Clear-Bold
